# Changing the Suite and module name for watchlist
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Rename "E Suite" -> "Watchlist" (row 6, column A)
$ws.Range("A6").Value = "Watchlist"

# Flip Runmode column from "Y" to "N" for every suite except the Watchlist one (row 6)
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C7").Value = "N"
# C6 (Watchlist row) keeps its original "Y" value

# Update the active/selected cell from C4 to B5
$ws.Range("B5").Select()
